# Append a new data row (row 35) to the worksheet, mirroring the existing
# Adafruit IO export rows (Timestamp, Feed Key, Value, Latitude, Longitude,
# Elevation). All columns in this sheet are stored as text, including the
# numeric-looking "Value" column, so we force column C to Text before
# writing it (otherwise Excel auto-coerces "25" into the number 25), then
# clear the temporary number-format override so the cell keeps the sheet's
# default (unstyled) formatting, matching the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 35

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "25"
$ws.Cells.Item($row, 3).ClearFormats()

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
